$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.391.54'
$ws.Range('E2').Value = '  +1.48%  '

$ws.Range('D3').Value = '1.907.37'
$ws.Range('E3').Value = '  +0.14%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('E5').Value = '  -1.68%  '

$ws.Range('E6').Value = '  -0.01%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4671'
$ws.Range('E7').Value = '  +0.63%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.4082'
$ws.Range('E8').Value = '  +0.59%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '47.75'
$ws.Range('E9').Value = '  -0.42%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.08014'
$ws.Range('E10').Value = '  +0.17%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '1.006'
$ws.Range('E11').Value = '  +0.49%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '22.29'
$ws.Range('E12').Value = '  +2.96%  '

$ws.Range('D13').Value = '1.920.68'
$ws.Range('E13').Value = '  +0.57%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.938'
$ws.Range('E14').Value = '  +0.39%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.121'
$ws.Range('E15').Value = '  +0.89%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '89.16'
$ws.Range('E16').Value = '  +0.30%  '

$ws.Range('E17').Value = '  +0.00%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.06603'
$ws.Range('E18').Value = '  +0.69%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.00001028'
$ws.Range('E19').Value = '  -0.56%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '17.71'
$ws.Range('E20').Value = '  +1.83%  '

$ws.Range('E21').Value = '  -0.15%  '

$ws.Range('D22').Value = '29.393.85'
$ws.Range('E22').Value = '  +1.44%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.529'
$ws.Range('E23').Value = '  +1.39%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.52'
$ws.Range('E24').Value = '  +3.12%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.212'
$ws.Range('E25').Value = '  -1.15%  '

$ws.Range('D26').Value = '2.118.99'

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '153.52'
$ws.Range('E27').Value = '  -2.63%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '19.76'
$ws.Range('E28').Value = '  +0.30%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.128'
$ws.Range('E29').Value = '  +1.59%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '5.702'
$ws.Range('E30').Value = '  +5.79%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '116.82'
$ws.Range('E31').Value = '  -1.63%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.073'
$ws.Range('E32').Value = '  +9.52%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.09476'

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.416'
$ws.Range('E34').Value = '  +0.13%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.576'
$ws.Range('E35').Value = '  -0.73%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.378'
$ws.Range('E36').Value = '  +1.71%  '

$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.06076'
$ws.Range('E37').Value = '  +0.04%  '

$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.02250'
$ws.Range('E38').Value = '  +1.26%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '8.360'
$ws.Range('E39').Value = '  -0.42%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.172'
$ws.Range('E40').Value = '  +0.86%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.5865'
$ws.Range('E41').Value = '  +1.33%  '

$ws.Range('E42').Value = '  +0.74%  '

$ws.Range('E43').Value = '  -0.05%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.300'
$ws.Range('E44').Value = '  +2.65%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.07729'
$ws.Range('E45').Value = '  +10.13%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.377'
$ws.Range('E46').Value = '  +2.37%  '

$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '12.19'
$ws.Range('E47').Value = '  +1.25%  '

$ws.Range('B48').Value = 'Decentraland'
$ws.Range('C48').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.5542'
$ws.Range('E48').Value = '  +0.99%  '

$ws.Range('E49').Value = '  +1.29%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '113.15'
$ws.Range('E50').Value = '  +0.89%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.2939'
$ws.Range('E51').Value = '  +5.71%  '
